# Atualização de bases das ligas, do dia: 28-04-2024 às 15:37
# Rows 6/7, 16/17, 86/87 had their match-data swapped (ordering fix),
# and the shared team names "SV Altldersdorf" / "FV Preussen Eberswalde"
# were corrected on the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Verbandsliga")

$ws.Cells.Item(6, 2).Value = 6781286
$ws.Cells.Item(6, 5).Value = "FV Preussen Eberswalde"
$ws.Cells.Item(6, 6).Value = "SV 1908 GW Ahrensfelde"
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = 5
$ws.Cells.Item(6, 9).Value = "A"
$ws.Cells.Item(6, 10).Value = 2.25
$ws.Cells.Item(6, 11).Value = 3.75
$ws.Cells.Item(6, 12).Value = 2.5
$ws.Cells.Item(6, 13).Value = 2.25
$ws.Cells.Item(6, 14).Value = 3.75
$ws.Cells.Item(6, 15).Value = 2.5
$ws.Cells.Item(6, 16).Value = -0.25
$ws.Cells.Item(6, 17).Value = 2.05
$ws.Cells.Item(6, 18).Value = 1.75
$ws.Cells.Item(6, 19).Value = 3.5
$ws.Cells.Item(6, 20).Value = 1.975
$ws.Cells.Item(6, 21).Value = 1.825
$ws.Cells.Item(6, 22).Value = -1
$ws.Cells.Item(6, 24).Value = 1.5
$ws.Cells.Item(6, 25).Value = -1
$ws.Cells.Item(6, 26).Value = 0.75
$ws.Cells.Item(6, 27).Value = 0.9750000000000001
$ws.Cells.Item(7, 2).Value = 6781300
$ws.Cells.Item(7, 5).Value = "SV Altldersdorf"
$ws.Cells.Item(7, 6).Value = "SV Frankonia Wernsdorf"
$ws.Cells.Item(7, 7).Value = 8
$ws.Cells.Item(7, 8).Value = 2
$ws.Cells.Item(7, 9).Value = "H"
$ws.Cells.Item(7, 10).Value = 2.2
$ws.Cells.Item(7, 11).Value = 3.5
$ws.Cells.Item(7, 12).Value = 2.7
$ws.Cells.Item(7, 13).Value = 1.727
$ws.Cells.Item(7, 14).Value = 4
$ws.Cells.Item(7, 15).Value = 3.5
$ws.Cells.Item(7, 16).Value = -0.5
$ws.Cells.Item(7, 17).Value = 1.775
$ws.Cells.Item(7, 18).Value = 2.025
$ws.Cells.Item(7, 19).Value = 3.25
$ws.Cells.Item(7, 20).Value = 1.925
$ws.Cells.Item(7, 21).Value = 1.875
$ws.Cells.Item(7, 22).Value = 0.7270000000000001
$ws.Cells.Item(7, 24).Value = -1
$ws.Cells.Item(7, 25).Value = 0.7749999999999999
$ws.Cells.Item(7, 26).Value = -1
$ws.Cells.Item(7, 27).Value = 0.925
$ws.Cells.Item(16, 2).Value = 7138607
$ws.Cells.Item(16, 5).Value = "Rot Weiss Walldorf II"
$ws.Cells.Item(16, 6).Value = "Turnerschaft OberRoden"
$ws.Cells.Item(16, 7).Value = 3
$ws.Cells.Item(16, 8).Value = 2
$ws.Cells.Item(16, 9).Value = "H"
$ws.Cells.Item(16, 10).Value = 2.25
$ws.Cells.Item(16, 11).Value = 3.75
$ws.Cells.Item(16, 12).Value = 2.5
$ws.Cells.Item(16, 13).Value = 2.25
$ws.Cells.Item(16, 14).Value = 3.8
$ws.Cells.Item(16, 15).Value = 2.45
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 1.8
$ws.Cells.Item(16, 18).Value = 2
$ws.Cells.Item(16, 19).Value = 3.75
$ws.Cells.Item(16, 20).Value = 1.95
$ws.Cells.Item(16, 21).Value = 1.85
$ws.Cells.Item(16, 22).Value = 1.25
$ws.Cells.Item(16, 23).Value = -1
$ws.Cells.Item(16, 25).Value = 0.8
$ws.Cells.Item(16, 26).Value = -1
$ws.Cells.Item(16, 27).Value = 0.95
$ws.Cells.Item(16, 28).Value = -1
$ws.Cells.Item(17, 2).Value = 7138608
$ws.Cells.Item(17, 5).Value = "SV UnterFlockenbach"
$ws.Cells.Item(17, 6).Value = "SC Dortelweil"
$ws.Cells.Item(17, 7).Value = 1
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(17, 9).Value = "D"
$ws.Cells.Item(17, 10).Value = 1.083
$ws.Cells.Item(17, 11).Value = 9
$ws.Cells.Item(17, 12).Value = 16
$ws.Cells.Item(17, 13).Value = 1.125
$ws.Cells.Item(17, 14).Value = 7.5
$ws.Cells.Item(17, 15).Value = 13
$ws.Cells.Item(17, 16).Value = -2.5
$ws.Cells.Item(17, 17).Value = 1.775
$ws.Cells.Item(17, 18).Value = 1.925
$ws.Cells.Item(17, 19).Value = 4.25
$ws.Cells.Item(17, 20).Value = 1.975
$ws.Cells.Item(17, 21).Value = 1.825
$ws.Cells.Item(17, 22).Value = -1
$ws.Cells.Item(17, 23).Value = 6.5
$ws.Cells.Item(17, 25).Value = -1
$ws.Cells.Item(17, 26).Value = 0.925
$ws.Cells.Item(17, 27).Value = -1
$ws.Cells.Item(17, 28).Value = 0.825
$ws.Cells.Item(86, 2).Value = 7511958
$ws.Cells.Item(86, 5).Value = "SpVgg EGC Wirges"
$ws.Cells.Item(86, 6).Value = "SG 2000 MulheimKarlich"
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(86, 8).Value = 1
$ws.Cells.Item(86, 10).Value = 4.333
$ws.Cells.Item(86, 11).Value = 4
$ws.Cells.Item(86, 12).Value = 1.571
$ws.Cells.Item(86, 13).Value = 4.2
$ws.Cells.Item(86, 15).Value = 1.571
$ws.Cells.Item(86, 16).Value = 1
$ws.Cells.Item(86, 17).Value = 1.875
$ws.Cells.Item(86, 18).Value = 1.925
$ws.Cells.Item(86, 19).Value = 3.75
$ws.Cells.Item(86, 20).Value = 1.925
$ws.Cells.Item(86, 21).Value = 1.875
$ws.Cells.Item(86, 22).Value = 3.2
$ws.Cells.Item(86, 25).Value = 0.875
$ws.Cells.Item(86, 27).Value = -1
$ws.Cells.Item(86, 28).Value = 0.875
$ws.Cells.Item(87, 2).Value = 7511976
$ws.Cells.Item(87, 5).Value = "DJK Bad Homburg"
$ws.Cells.Item(87, 6).Value = "SG Bornheim 1945 GrunWeiss"
$ws.Cells.Item(87, 7).Value = 4
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 2
$ws.Cells.Item(87, 11).Value = 3.75
$ws.Cells.Item(87, 12).Value = 2.9
$ws.Cells.Item(87, 13).Value = 1.8
$ws.Cells.Item(87, 15).Value = 3.3
$ws.Cells.Item(87, 16).Value = -0.5
$ws.Cells.Item(87, 17).Value = 1.85
$ws.Cells.Item(87, 18).Value = 1.95
$ws.Cells.Item(87, 19).Value = 3.5
$ws.Cells.Item(87, 20).Value = 1.975
$ws.Cells.Item(87, 21).Value = 1.825
$ws.Cells.Item(87, 22).Value = 0.8
$ws.Cells.Item(87, 25).Value = 0.8500000000000001
$ws.Cells.Item(87, 27).Value = 0.9750000000000001
$ws.Cells.Item(87, 28).Value = -1
